$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 42496
$ws.Range("E2").Value = 1206
$ws.Range("F2").Value = 1206
$ws.Range("G2").Value = 2111
$ws.Range("H2").Value = 1670
$ws.Range("I2").Value = 1644
$ws.Range("J2").Value = 26
$ws.Range("K2").Value = 204257
$ws.Range("L2").Value = 186198
$ws.Range("M2").Value = 18060
$ws.Range("N2").Value = 17714
$ws.Range("O2").Value = 345
$ws.Range("P2").Value = 5379
$ws.Range("Q2").Value = 6968
$ws.Range("R2").Value = -4160
$ws.Range("S2").Value = -2181
$ws.Range("T2").Value = 79
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 2.84
$ws.Range("X2").Value = 3.93
$ws.Range("Y2").Value = 10.74
$ws.Range("Z2").Value = 0.87
$ws.Range("AA2").Value = 1031.02
$ws.Range("AB2").Value = 247.07
$ws.Range("AC2").Value = 1529
$ws.Range("AD2").Value = 7.03
$ws.Range("AE2").Value = 17355
$ws.Range("AF2").Value = 0.62
$ws.Range("AG2").Value = 550
$ws.Range("AH2").Value = 5.12
$ws.Range("AI2").Value = 34.14
$ws.Range("AJ2").Value = 107572390
$ws.Range("U2").ClearContents()

# Row 3
$ws.Range("D3").Value = 46985
$ws.Range("E3").Value = 2122
$ws.Range("F3").Value = 2122
$ws.Range("G3").Value = 2057
$ws.Range("H3").Value = 1605
$ws.Range("I3").Value = 1579
$ws.Range("J3").Value = 26
$ws.Range("K3").Value = 226209
$ws.Range("L3").Value = 206282
$ws.Range("M3").Value = 19927
$ws.Range("N3").Value = 19651
$ws.Range("O3").Value = 276
$ws.Range("P3").Value = 5379
$ws.Range("Q3").Value = 9221
$ws.Range("R3").Value = -9586
$ws.Range("S3").Value = -565
$ws.Range("T3").Value = 46
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 4.52
$ws.Range("X3").Value = 3.42
$ws.Range("Y3").Value = 8.59
$ws.Range("Z3").Value = 0.75
$ws.Range("AA3").Value = 1035.18
$ws.Range("AB3").Value = 281.79
$ws.Range("AC3").Value = 1468
$ws.Range("AD3").Value = 7.94
$ws.Range("AE3").Value = 19252
$ws.Range("AF3").Value = 0.61
$ws.Range("AG3").Value = 620
$ws.Range("AH3").Value = 5.32
$ws.Range("AI3").Value = 40.08
$ws.Range("AJ3").Value = 107572390
$ws.Range("U3").ClearContents()

# Row 4
$ws.Range("D4").Value = 74295
$ws.Range("E4").Value = -298
$ws.Range("F4").Value = -298
$ws.Range("G4").Value = -321
$ws.Range("H4").Value = 148
$ws.Range("I4").Value = 120
$ws.Range("J4").Value = 28
$ws.Range("K4").Value = 267208
$ws.Range("L4").Value = 248851
$ws.Range("M4").Value = 18357
$ws.Range("N4").Value = 18129
$ws.Range("O4").Value = 229
$ws.Range("P4").Value = 5379
$ws.Range("Q4").Value = 40166
$ws.Range("R4").Value = -39433
$ws.Range("S4").Value = -641
$ws.Range("T4").Value = 55
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = -0.4
$ws.Range("X4").Value = 0.2
$ws.Range("Y4").Value = 0.78
$ws.Range("Z4").Value = 0.06
$ws.Range("AA4").Value = 1355.61
$ws.Range("AB4").Value = 252.6
$ws.Range("AC4").Value = 112
$ws.Range("AD4").Value = 113.87
$ws.Range("AE4").Value = 17760
$ws.Range("AF4").Value = 0.72
$ws.Range("AG4").Value = 200
$ws.Range("AH4").Value = 1.57
$ws.Range("AI4").Value = 170.16
$ws.Range("AJ4").Value = 107572390
$ws.Range("U4").ClearContents()

# Row 5
$ws.Range("D5").Value = 71397
$ws.Range("E5").Value = 2466
$ws.Range("F5").Value = 2466
$ws.Range("G5").Value = 2519
$ws.Range("H5").Value = 1928
$ws.Range("I5").Value = 1900
$ws.Range("J5").Value = 28
$ws.Range("K5").Value = 303439
$ws.Range("L5").Value = 278915
$ws.Range("M5").Value = 24524
$ws.Range("N5").Value = 24163
$ws.Range("O5").Value = 361
$ws.Range("P5").Value = 8068
$ws.Range("Q5").Value = 29845
$ws.Range("R5").Value = -35261
$ws.Range("S5").Value = 5056
$ws.Range("T5").Value = 60
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 3.45
$ws.Range("X5").Value = 2.7
$ws.Range("Y5").Value = 9.119999999999999
$ws.Range("Z5").Value = 0.68
$ws.Range("AA5").Value = 1137.31
$ws.Range("AB5").Value = 211.51
$ws.Range("AC5").Value = 1277
$ws.Range("AD5").Value = 5.95
$ws.Range("AE5").Value = 15503
$ws.Range("AF5").Value = 0.49
$ws.Range("AG5").Value = 360
$ws.Range("AH5").Value = 4.74
$ws.Range("AI5").Value = 29.53
$ws.Range("AJ5").Value = 161358585
$ws.Range("U5").ClearContents()

# Row 6
$ws.Range("D6").Value = 57869
$ws.Range("E6").Value = 726
$ws.Range("F6").Value = 726
$ws.Range("G6").Value = 764
$ws.Range("H6").Value = 566
$ws.Range("I6").Value = 548
$ws.Range("K6").Value = 319207
$ws.Range("L6").Value = 296234
$ws.Range("M6").Value = 22973
$ws.Range("N6").Value = 22709
$ws.Range("P6").Value = 8068
$ws.Range("Q6").Value = 16692
$ws.Range("R6").Value = -16748
$ws.Range("S6").Value = 435
$ws.Range("T6").Value = 37
$ws.Range("V6").Value = 1000
$ws.Range("W6").Value = 1.25
$ws.Range("X6").Value = 0.98
$ws.Range("Y6").Value = 2.42
$ws.Range("Z6").Value = 0.18
$ws.Range("AA6").Value = 1289.49
$ws.Range("AB6").Value = 192.28
$ws.Range("AC6").Value = 339
$ws.Range("AD6").Value = 14.05
$ws.Range("AE6").Value = 14570
$ws.Range("AF6").Value = 0.33
$ws.Range("AG6").Value = 100
$ws.Range("AH6").Value = 2.1
$ws.Range("AI6").Value = 28.46
$ws.Range("AJ6").Value = 161358585
$ws.Range("U6").ClearContents()

# Row 7
$ws.Range("D7").Value = 47430
$ws.Range("G7").Value = 1980
$ws.Range("H7").Value = 1500
$ws.Range("I7").Value = 1480
$ws.Range("K7").Value = 341490
$ws.Range("L7").Value = 314140
$ws.Range("M7").Value = 27350
$ws.Range("N7").Value = 27120
$ws.Range("X7").Value = 3.16
$ws.Range("Y7").Value = 5.94
$ws.Range("Z7").Value = 0.45
$ws.Range("AA7").Value = 1148.59
$ws.Range("AC7").Value = 917
$ws.Range("AD7").Value = 3.92
$ws.Range("AE7").Value = 17400
$ws.Range("AF7").Value = 0.21
$ws.Range("AG7").Value = 232
$ws.Range("AH7").Value = 6.45
$ws.Range("AI7").Value = 25.28
$ws.Range("E7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()

# Row 8
$ws.Range("D8").Value = 49770
$ws.Range("G8").Value = 1370
$ws.Range("H8").Value = 1040
$ws.Range("I8").Value = 1040
$ws.Range("K8").Value = 358180
$ws.Range("L8").Value = 330170
$ws.Range("M8").Value = 28020
$ws.Range("N8").Value = 27780
$ws.Range("X8").Value = 2.09
$ws.Range("Y8").Value = 3.79
$ws.Range("Z8").Value = 0.3
$ws.Range("AA8").Value = 1178.34
$ws.Range("AC8").Value = 645
$ws.Range("AD8").Value = 5.58
$ws.Range("AE8").Value = 17824
$ws.Range("AF8").Value = 0.2
$ws.Range("AG8").Value = 209
$ws.Range("AH8").Value = 5.81
$ws.Range("AI8").Value = 32.38
$ws.Range("E8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()

# Row 9
$ws.Range("G9").Value = 1560
$ws.Range("H9").Value = 1180
$ws.Range("I9").Value = 1180
$ws.Range("K9").Value = 373730
$ws.Range("L9").Value = 344910
$ws.Range("M9").Value = 28830
$ws.Range("N9").Value = 28590
$ws.Range("Y9").Value = 4.19
$ws.Range("Z9").Value = 0.32
$ws.Range("AA9").Value = 1196.36
$ws.Range("AC9").Value = 731
$ws.Range("AD9").Value = 4.92
$ws.Range("AE9").Value = 18344
$ws.Range("AF9").Value = 0.2
$ws.Range("AG9").Value = 192
$ws.Range("AH9").Value = 26.28
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("AI9").ClearContents()
